$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 6945
$ws1.Range("F6").Value = 1855
$ws1.Range("F7").Value = 6435
$ws1.Range("F9").Value = 1996
$ws1.Range("F10").Value = 537
$ws1.Range("F11").Value = 30
$ws1.Range("F17").Value = 8242
$ws1.Range("F22").Value = 1783
$ws1.Range("F23").Value = 855
$ws1.Range("F30").Value = 1924
$ws1.Range("F31").Value = 829
$ws1.Range("F32").Value = 431
$ws1.Range("F34").Value = 10
$ws1.Range("F35").Value = 145
$ws1.Range("F36").Value = 123
$ws1.Range("F37").Value = 94
$ws1.Range("F38").Value = 3943

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 292

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 6945
$ws4.Range("F10").Value = 1855
$ws4.Range("F11").Value = 6435
$ws4.Range("F13").Value = 1996
$ws4.Range("F15").Value = 537
$ws4.Range("F23").Value = 8242
$ws4.Range("F28").Value = 1783
$ws4.Range("F29").Value = 855
$ws4.Range("F33").Value = 1924
$ws4.Range("F34").Value = 829
$ws4.Range("F36").Value = 431
$ws4.Range("F39").Value = 10
$ws4.Range("F41").Value = 123
$ws4.Range("F42").Value = 94
$ws4.Range("F44").Value = 3943
